$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 30 de Septiembre de 2020 a las 19:57'
$ws.Cells.Item(4, 2).Value = 7421094
$ws.Cells.Item(4, 3).Value = 14948
$ws.Cells.Item(4, 4).Value = 4670556
$ws.Cells.Item(4, 5).Value = 2539322
$ws.Cells.Item(4, 7).Value = 431
$ws.Cells.Item(4, 8).Value = 211216
$ws.Cells.Item(5, 2).Value = 6305643
$ws.Cells.Item(5, 3).Value = 82124
$ws.Cells.Item(5, 4).Value = 5263685
$ws.Cells.Item(5, 5).Value = 943342
$ws.Cells.Item(5, 7).Value = 1087
$ws.Cells.Item(5, 8).Value = 98616
$ws.Cells.Item(10, 2).Value = 769188
$ws.Cells.Item(10, 3).Value = 11016
$ws.Cells.Item(10, 7).Value = 177
$ws.Cells.Item(10, 8).Value = 31791
$ws.Cells.Item(14, 2).Value = 563535
$ws.Cells.Item(14, 3).Value = 12845
$ws.Cells.Item(14, 5).Value = 435252
$ws.Cells.Item(14, 7).Value = 63
$ws.Cells.Item(14, 8).Value = 31956
$ws.Cells.Item(21, 2).Value = 318663
$ws.Cells.Item(21, 3).Value = 1391
$ws.Cells.Item(21, 4).Value = 279749
$ws.Cells.Item(21, 5).Value = 30719
$ws.Cells.Item(21, 7).Value = 65
$ws.Cells.Item(21, 8).Value = 8195
$ws.Cells.Item(25, 2).Value = 291929
$ws.Cells.Item(25, 3).Value = 1463
$ws.Cells.Item(25, 5).Value = 26369
$ws.Cells.Item(25, 7).Value = 4
$ws.Cells.Item(25, 8).Value = 9560
$ws.Cells.Item(27, 2).Value = 243895
$ws.Cells.Item(27, 3).Value = 6969
$ws.Cells.Item(27, 4).Value = 174232
$ws.Cells.Item(27, 5).Value = 68111
$ws.Cells.Item(27, 7).Value = 24
$ws.Cells.Item(27, 8).Value = 1552
$ws.Cells.Item(34, 2).Value = 123653
$ws.Cells.Item(34, 3).Value = 2470
$ws.Cells.Item(34, 4).Value = 102715
$ws.Cells.Item(34, 5).Value = 18744
$ws.Cells.Item(34, 7).Value = 42
$ws.Cells.Item(34, 8).Value = 2194
$ws.Cells.Item(43, 2).Value = 94190
$ws.Cells.Item(43, 3).Value = 1100
$ws.Cells.Item(43, 4).Value = 83724
$ws.Cells.Item(43, 5).Value = 10047
$ws.Cells.Item(43, 7).Value = 3
$ws.Cells.Item(43, 8).Value = 419
$ws.Cells.Item(53, 1).Value = 'Etiopia'
$ws.Cells.Item(53, 2).Value = 75368
$ws.Cells.Item(53, 3).Value = 784
$ws.Cells.Item(53, 4).Value = 31204
$ws.Cells.Item(53, 5).Value = 42966
$ws.Cells.Item(53, 7).Value = 7
$ws.Cells.Item(53, 8).Value = 1198
$ws.Cells.Item(54, 1).Value = 'Costa Rica'
$ws.Cells.Item(54, 2).Value = 74604
$ws.Cells.Item(54, 4).Value = 30703
$ws.Cells.Item(54, 5).Value = 43021
$ws.Cells.Item(54, 8).Value = 880
$ws.Cells.Item(71, 1).Value = 'Libano'
$ws.Cells.Item(71, 2).Value = 39634
$ws.Cells.Item(71, 3).Value = 1257
$ws.Cells.Item(71, 4).Value = 17565
$ws.Cells.Item(71, 5).Value = 21702
$ws.Cells.Item(71, 7).Value = 6
$ws.Cells.Item(71, 8).Value = 367
$ws.Cells.Item(72, 1).Value = 'Afganistan'
$ws.Cells.Item(72, 2).Value = 39268
$ws.Cells.Item(72, 3).Value = 14
$ws.Cells.Item(72, 4).Value = 32789
$ws.Cells.Item(72, 5).Value = 5021
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 1458
$ws.Cells.Item(73, 1).Value = 'Kenia'
$ws.Cells.Item(73, 2).Value = 38529
$ws.Cells.Item(73, 3).Value = 151
$ws.Cells.Item(73, 4).Value = 24908
$ws.Cells.Item(73, 5).Value = 12910
$ws.Cells.Item(73, 7).Value = 4
$ws.Cells.Item(73, 8).Value = 711
$ws.Cells.Item(74, 2).Value = 36155
$ws.Cells.Item(74, 3).Value = 415
$ws.Cells.Item(74, 5).Value = 10987
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 1804
$ws.Cells.Item(87, 1).Value = 'Tunez'
$ws.Cells.Item(87, 2).Value = 18413
$ws.Cells.Item(87, 3).Value = 1008
$ws.Cells.Item(87, 4).Value = 5032
$ws.Cells.Item(87, 5).Value = 13116
$ws.Cells.Item(87, 7).Value = 19
$ws.Cells.Item(87, 8).Value = 265
$ws.Cells.Item(88, 1).Value = 'Republica de Macedonia'
$ws.Cells.Item(88, 2).Value = 17977
$ws.Cells.Item(88, 3).Value = 191
$ws.Cells.Item(88, 4).Value = 14959
$ws.Cells.Item(88, 5).Value = 2279
$ws.Cells.Item(88, 7).Value = 2
$ws.Cells.Item(88, 8).Value = 739
$ws.Cells.Item(110, 2).Value = 8728
$ws.Cells.Item(110, 3).Value = 172
$ws.Cells.Item(110, 4).Value = 5232
$ws.Cells.Item(110, 5).Value = 3435
$ws.Cells.Item(110, 7).Value = 2
$ws.Cells.Item(110, 8).Value = 61
$ws.Cells.Item(121, 2).Value = 5482
$ws.Cells.Item(121, 3).Value = 20
$ws.Cells.Item(121, 4).Value = 4912
$ws.Cells.Item(121, 5).Value = 461
$ws.Cells.Item(121, 7).Value = 1
$ws.Cells.Item(121, 8).Value = 109
$ws.Cells.Item(145, 2).Value = 3118
$ws.Cells.Item(145, 3).Value = 17
$ws.Cells.Item(145, 4).Value = 2453
$ws.Cells.Item(145, 5).Value = 534
$ws.Cells.Item(146, 5).Value = 461
$ws.Cells.Item(146, 7).Value = 1
$ws.Cells.Item(146, 8).Value = 35
$ws.Cells.Item(149, 2).Value = 2704
$ws.Cells.Item(149, 3).Value = 4
$ws.Cells.Item(149, 5).Value = 1365
$ws.Cells.Item(153, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(153, 2).Value = 2050
$ws.Cells.Item(153, 3).Value = 84
$ws.Cells.Item(153, 4).Value = 1432
$ws.Cells.Item(153, 5).Value = 565
$ws.Cells.Item(153, 8).Value = 53
$ws.Cells.Item(154, 1).Value = 'Yemen'
$ws.Cells.Item(154, 2).Value = 2034
$ws.Cells.Item(154, 3).Value = 3
$ws.Cells.Item(154, 4).Value = 1286
$ws.Cells.Item(154, 5).Value = 161
$ws.Cells.Item(154, 8).Value = 587
$ws.Cells.Item(155, 1).Value = 'Uruguay'
$ws.Cells.Item(155, 2).Value = 2033
$ws.Cells.Item(155, 4).Value = 1771
$ws.Cells.Item(155, 5).Value = 214
$ws.Cells.Item(155, 8).Value = 48
$ws.Cells.Item(156, 1).Value = 'Burkina Faso'
$ws.Cells.Item(156, 2).Value = 2032
$ws.Cells.Item(156, 4).Value = 1309
$ws.Cells.Item(156, 5).Value = 665
$ws.Cells.Item(156, 8).Value = 58
$ws.Cells.Item(207, 1).Value = 'Santa Lucia'
$ws.Cells.Item(208, 1).Value = 'Nueva Caledonia'
